# "about to click go" - sync pass: bumps _updated/_created timestamps to the
# new sync time (4/29/2023, 10:51:20 PM), mints fresh _uid values for every
# touched record, regenerates a couple of short ids, and replaces the single
# "Select Option Test" tag-def option with two new options (Select Option 1 /
# Select Option 2), wiring up their Tags rows.

$wb = $excel.ActiveWorkbook

# The app's JS `Date.toLocaleString()` renders the time/AM-PM separator as a
# NARROW NO-BREAK SPACE (U+202F), not an ordinary space - match that so the
# new timestamp strings are byte-identical in shape to the old ones.
$nbsp          = [char]0x202F
$newStamp      = "4/29/2023, 10:51:20" + $nbsp + "PM"
$newStampCDT   = "4/29/2023, 10:51:20" + $nbsp + "PM CDT"
$newStampIso   = "2023-04-29T22:51:20"

# ---------------------------------------------------------------------
# Overview
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B2").Value = $newStampCDT
$ws.Range("B8").Value = 5          # tagDefs count active
$ws.Range("B9").Value = 4          # tags count active

# ---------------------------------------------------------------------
# Defs
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Defs")
$ws.Range("A2").Value = "lh2vk0q0-0db6"
$ws.Range("B2").Value = $newStamp
$ws.Range("C2").Value = $newStamp

$ws.Range("A3").Value = "lh2vk0q3-0dtj"
$ws.Range("B3").Value = $newStamp
$ws.Range("C3").Value = $newStamp

$ws.Range("A4").Value = "lh2vk0q4-0a06"
$ws.Range("C4").Value = $newStamp

# ---------------------------------------------------------------------
# Point Defs
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Point Defs")
$ws.Range("A2").Value = "lh2vk0q2-ojke"
$ws.Range("B2").Value = $newStamp
$ws.Range("C2").Value = $newStamp

$ws.Range("A3").Value = "lh2vk0q4-0o48"
$ws.Range("B3").Value = $newStamp
$ws.Range("C3").Value = $newStamp
$ws.Range("F3").Value = "0iu2"

$ws.Range("A4").Value = "lh2vk0q5-klza"
$ws.Range("B4").Value = $newStamp
$ws.Range("C4").Value = $newStamp

$ws.Range("A5").Value = "lh2vk0q5-1wpj"
$ws.Range("B5").Value = $newStamp
$ws.Range("C5").Value = $newStamp

# ---------------------------------------------------------------------
# Entry
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Entry")
$ws.Range("A2").Value = "lh2vk0q6-0m06"
$ws.Range("B2").Value = $newStamp
$ws.Range("C2").Value = $newStamp

$ws.Range("A3").Value = "lh2vk0q7-0kef"
$ws.Range("B3").Value = $newStamp
$ws.Range("C3").Value = $newStamp
$ws.Range("F3").Value = "lh2vk0q8-3am6"
$ws.Range("G3").Value = $newStampIso

# ---------------------------------------------------------------------
# Entry Points
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Entry Points")
$ws.Range("A2").Value = "lh2vk0q6-idb9"
$ws.Range("B2").Value = $newStamp
$ws.Range("C2").Value = $newStamp

$ws.Range("A3").Value = "lh2vk0q7-r0s5"
$ws.Range("B3").Value = $newStamp
$ws.Range("C3").Value = $newStamp

# ---------------------------------------------------------------------
# Tag Defs  (row 4 "Select Option Test" -> "Select Option To Delete",
#            plus two brand-new option rows 5 & 6)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Tag Defs")
$ws.Range("A2").Value = "lh2vk0q9-0iha"
$ws.Range("B2").Value = $newStamp
$ws.Range("C2").Value = $newStamp
$ws.Range("E2").Value = "0x2q"

$ws.Range("A3").Value = "lh2vk0q9-5k8p"
$ws.Range("B3").Value = $newStamp
$ws.Range("C3").Value = $newStamp

$ws.Range("A4").Value = "lh2vk0qa-qbib"
$ws.Range("B4").Value = $newStamp
$ws.Range("C4").Value = $newStamp
$ws.Range("F4").Value = "Select Option To Delete"

$ws.Range("A5").Value = "lh2vk0qa-4x8q"
$ws.Range("B5").Value = $newStamp
$ws.Range("C5").Value = $newStamp
$ws.Range("D5").Value = "'FALSE"
$ws.Range("E5").Value = "0vva"
$ws.Range("F5").Value = "Select Option 1"

$ws.Range("A6").Value = "lh2vk0qb-gn9e"
$ws.Range("B6").Value = $newStamp
$ws.Range("C6").Value = $newStamp
$ws.Range("D6").Value = "'FALSE"
$ws.Range("E6").Value = "0vvb"
$ws.Range("F6").Value = "Select Option 2"

# ---------------------------------------------------------------------
# Tags  (two brand-new tag rows 4 & 5 pointing at the new options)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Tags")
$ws.Range("A2").Value = "lh2vk0qb-bq3o"
$ws.Range("B2").Value = $newStamp
$ws.Range("C2").Value = $newStamp

$ws.Range("A3").Value = "lh2vk0qc-98am"
$ws.Range("B3").Value = $newStamp
$ws.Range("C3").Value = $newStamp

$ws.Range("A4").Value = "lh2vk0qc-qss9"
$ws.Range("B4").Value = $newStamp
$ws.Range("C4").Value = $newStamp
$ws.Range("D4").Value = "'FALSE"
$ws.Range("E4").Value = "0m7w"
$ws.Range("F4").Value = "8esq"
$ws.Range("G4").Value = "0vva"

$ws.Range("A5").Value = "lh2vk0qd-95rs"
$ws.Range("B5").Value = $newStamp
$ws.Range("C5").Value = $newStamp
$ws.Range("D5").Value = "'FALSE"
$ws.Range("E5").Value = "0m7w"
$ws.Range("F5").Value = "8esq"
$ws.Range("G5").Value = "0vvb"
